# Apply updated "dSF" (column F) values as part of a data repull / mean recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F8").Value = -15
$ws.Range("F9").Value = 5
$ws.Range("F13").Value = -9
$ws.Range("F16").Value = 3
